$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3795.625
$ws.Range("I98").Value = 3868
$ws.Range("J98").Value = 2999.5
$ws.Range("K98").Value = 3868
$ws.Range("L98").Value = 2999.5
$ws.Range("M98").Value = -2370
$ws.Range("N98").Value = -5995.5
$ws.Range("H122").Value = 3795.625
$ws.Range("I122").Value = 3868
$ws.Range("J122").Value = 2999.5
$ws.Range("K122").Value = 11604
$ws.Range("L122").Value = 8998.5
$ws.Range("M122").Value = -9154
$ws.Range("N122").Value = -13898.5
$ws.Range("H132").Value = 4043.149
$ws.Range("I132").Value = 4465.143
$ws.Range("J132").Value = 498.4
$ws.Range("K132").Value = 13395.429
$ws.Range("L132").Value = 1495.2
$ws.Range("M132").Value = -10865.429
$ws.Range("N132").Value = -6555.2
$ws.Range("H137").Value = 2943150
$ws.Range("I137").Value = 3335181.8
$ws.Range("K137").Value = 10005545.4
$ws.Range("M137").Value = -10002995.4
$ws.Range("H141").Value = 4058.9092
$ws.Range("I141").Value = 4082.25
$ws.Range("K141").Value = 12246.75
$ws.Range("M141").Value = -7066.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1422.421
$ws.Range("I2").Value = 1056.2858
$ws.Range("J2").Value = 2447.6
$ws.Range("K2").Value = 1056.2858
$ws.Range("L2").Value = 2447.6
$ws.Range("M2").Value = -943.2858000000001
$ws.Range("N2").Value = -2673.6
$ws.Range("H32").Value = 6863.8125
$ws.Range("I32").Value = 7188.067
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 7188.067
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -6901.067
$ws.Range("N32").Value = -2574
$ws.Range("H61").Value = 2913.2307
$ws.Range("I61").Value = 1980.3077
$ws.Range("J61").Value = 3846.1538
$ws.Range("K61").Value = 1980.3077
$ws.Range("L61").Value = 3846.1538
$ws.Range("M61").Value = -1768.3077
$ws.Range("N61").Value = -4270.1538
$ws.Range("H74").Value = 160977.58
$ws.Range("I74").Value = 200261.25
$ws.Range("K74").Value = 200261.25
$ws.Range("M74").Value = -199387.25
$ws.Range("H77").Value = 160977.58
$ws.Range("I77").Value = 200261.25
$ws.Range("K77").Value = 1001306.25
$ws.Range("M77").Value = -996938.25
$ws.Range("H116").Value = 1422.421
$ws.Range("I116").Value = 1056.2858
$ws.Range("J116").Value = 2447.6
$ws.Range("K116").Value = 1056.2858
$ws.Range("L116").Value = 2447.6
$ws.Range("M116").Value = 1237.7142
$ws.Range("N116").Value = -7035.6
$ws.Range("H124").Value = 66086.28999999999
$ws.Range("J124").Value = 66086.28999999999
$ws.Range("L124").Value = 66086.28999999999
$ws.Range("N124").Value = -75906.28999999999
$ws.Range("H132").Value = 2305.5454
$ws.Range("I132").Value = 2117.9167
$ws.Range("J132").Value = 3149.875
$ws.Range("K132").Value = 6353.750100000001
$ws.Range("L132").Value = 9449.625
$ws.Range("M132").Value = -3823.750100000001
$ws.Range("N132").Value = -14509.625
$ws.Range("H136").Value = 2913.2307
$ws.Range("I136").Value = 1980.3077
$ws.Range("J136").Value = 3846.1538
$ws.Range("K136").Value = 5940.9231
$ws.Range("L136").Value = 11538.4614
$ws.Range("M136").Value = -3390.9231
$ws.Range("N136").Value = -16638.4614
$ws.Range("H139").Value = 75486.25
$ws.Range("J139").Value = 75486.25
$ws.Range("L139").Value = 75486.25
$ws.Range("N139").Value = -85766.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1422.421
$ws.Range("I3").Value = 1056.2858
$ws.Range("J3").Value = 2447.6
$ws.Range("K3").Value = 1056.2858
$ws.Range("L3").Value = 2447.6
$ws.Range("M3").Value = -942.2858000000001
$ws.Range("N3").Value = -2675.6
$ws.Range("H94").Value = 86958620
$ws.Range("I94").Value = 133334984
$ws.Range("K94").Value = 133334984
$ws.Range("M94").Value = -133334533
$ws.Range("H105").Value = 9288621
$ws.Range("I105").Value = 479076.2
$ws.Range("J105").Value = 35717256
$ws.Range("K105").Value = 479076.2
$ws.Range("L105").Value = 35717256
$ws.Range("M105").Value = -477329.2
$ws.Range("N105").Value = -35720750
$ws.Range("H107").Value = 2565151.8
$ws.Range("I107").Value = 3497410.5
$ws.Range("J107").Value = 1440.125
$ws.Range("K107").Value = 3497410.5
$ws.Range("L107").Value = 1440.125
$ws.Range("M107").Value = -3495490.5
$ws.Range("N107").Value = -5280.125
$ws.Range("H134").Value = 2328.9473
$ws.Range("I134").Value = 2027.1086
$ws.Range("K134").Value = 6081.325800000001
$ws.Range("M134").Value = -3546.325800000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2477.4
$ws.Range("I16").Value = 2477.4
$ws.Range("K16").Value = 2477.4
$ws.Range("M16").Value = -2190.4
$ws.Range("H31").Value = 4634966
$ws.Range("I31").Value = 5599.5557
$ws.Range("J31").Value = 6949649.5
$ws.Range("K31").Value = 5599.5557
$ws.Range("L31").Value = 6949649.5
$ws.Range("M31").Value = -5304.5557
$ws.Range("N31").Value = -6950239.5
$ws.Range("H34").Value = 4634966
$ws.Range("I34").Value = 5599.5557
$ws.Range("J34").Value = 6949649.5
$ws.Range("K34").Value = 5599.5557
$ws.Range("L34").Value = 6949649.5
$ws.Range("M34").Value = -5397.5557
$ws.Range("N34").Value = -6950053.5
$ws.Range("H58").Value = 3174.4666
$ws.Range("I58").Value = 3249.5715
$ws.Range("K58").Value = 3249.5715
$ws.Range("M58").Value = -3046.5715
$ws.Range("H94").Value = 909.1818
$ws.Range("J94").Value = 1030.9286
$ws.Range("L94").Value = 1030.9286
$ws.Range("N94").Value = -1932.9286
$ws.Range("H99").Value = 1872.5
$ws.Range("I99").Value = 1750
$ws.Range("J99").Value = 1995
$ws.Range("K99").Value = 1750
$ws.Range("L99").Value = 1995
$ws.Range("M99").Value = -252
$ws.Range("N99").Value = -4991
$ws.Range("H105").Value = 1838.125
$ws.Range("I105").Value = 1386.4286
$ws.Range("K105").Value = 1386.4286
$ws.Range("H107").Value = 1852376.1
$ws.Range("I107").Value = 2500351.5
$ws.Range("K107").Value = 2500351.5
$ws.Range("M107").Value = -2498431.5
$ws.Range("H113").Value = 2477.4
$ws.Range("I113").Value = 2477.4
$ws.Range("K113").Value = 2477.4
$ws.Range("M113").Value = -307.4000000000001
$ws.Range("H126").Value = 1872.5
$ws.Range("I126").Value = 1750
$ws.Range("J126").Value = 1995
$ws.Range("K126").Value = 5250
$ws.Range("L126").Value = 5985
$ws.Range("M126").Value = -2780
$ws.Range("N126").Value = -10925
$ws.Range("H132").Value = 15153881
$ws.Range("I132").Value = 2155.125
$ws.Range("J132").Value = 23812010
$ws.Range("K132").Value = 6465.375
$ws.Range("L132").Value = 71436030
$ws.Range("M132").Value = -3935.375
$ws.Range("N132").Value = -71441090
$ws.Range("H136").Value = 3174.4666
$ws.Range("I136").Value = 3249.5715
$ws.Range("K136").Value = 9748.7145
$ws.Range("M136").Value = -7198.7145
$ws.Range("H141").Value = 244537.56
$ws.Range("J141").Value = 259475.08
$ws.Range("L141").Value = 259475.08
$ws.Range("N141").Value = -269835.08

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 31385640
$ws.Range("I70").Value = 35861450
$ws.Range("K70").Value = 35861450
$ws.Range("M70").Value = -35861180
$ws.Range("H73").Value = 31385640
$ws.Range("I73").Value = 35861450
$ws.Range("K73").Value = 35861450
$ws.Range("M73").Value = -35860514
$ws.Range("H97").Value = 1450.8667
$ws.Range("I97").Value = 1404.5
$ws.Range("J97").Value = 1543.6
$ws.Range("K97").Value = 1404.5
$ws.Range("L97").Value = 1543.6
$ws.Range("M97").Value = -908.5
$ws.Range("N97").Value = -2535.6
$ws.Range("H122").Value = 5548
$ws.Range("I122").Value = 3028
$ws.Range("K122").Value = 9084
$ws.Range("M122").Value = -6634
$ws.Range("H132").Value = 2183.3704
$ws.Range("I132").Value = 1848.0588
$ws.Range("K132").Value = 5544.1764
$ws.Range("M132").Value = -3014.1764

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6960.2
$ws.Range("I7").Value = 2763.9092
$ws.Range("K7").Value = 2763.9092
$ws.Range("M7").Value = -2651.9092
$ws.Range("H16").Value = 1337.3684
$ws.Range("I16").Value = 1387.2667
$ws.Range("K16").Value = 1387.2667
$ws.Range("M16").Value = -1217.2667
$ws.Range("H122").Value = 3954.8333
$ws.Range("J122").Value = 1943
$ws.Range("L122").Value = 5829
$ws.Range("N122").Value = -10729
$ws.Range("H126").Value = 6960.2
$ws.Range("I126").Value = 2763.9092
$ws.Range("K126").Value = 8291.7276
$ws.Range("M126").Value = -5821.7276
$ws.Range("H132").Value = 4018.1035
$ws.Range("J132").Value = 3651.111
$ws.Range("L132").Value = 10953.333
$ws.Range("N132").Value = -16013.333
$ws.Range("H136").Value = 3661.9048
$ws.Range("I136").Value = 3542.8572
$ws.Range("K136").Value = 10628.5716
$ws.Range("M136").Value = -8078.571599999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H122").Value = 15626949
$ws.Range("I122").Value = 2056.8462
$ws.Range("K122").Value = 6170.5386
$ws.Range("M122").Value = -3720.5386
$ws.Range("H132").Value = 3878.5625
$ws.Range("I132").Value = 4218.7144
$ws.Range("K132").Value = 12656.1432
$ws.Range("M132").Value = -10126.1432

Write-Host "Applied all changes"